$d = $word.ActiveDocument

# A <w:br/> manual line break renders as Chr(11) (vertical tab) in Range.Text.
$vt = [char]11

$findText = $vt + "Project Metadata:" + $vt + `
    "- Author: MuhammadAbdullah95 (ma2404374@gmail.com)" + $vt + `
    "- Python >= 3.11" + $vt + `
    "- Version: 0.1.0"
$replaceText = $vt

$range = $d.Content
$result = $range.Find.Execute($findText, $false, $false, $false, $false, $false, $true, 1, $false, $replaceText, 2)

if (-not $result) {
    throw "Find/Replace for Project Metadata block did not match anything."
}

Write-Output "Replace result: $result"
